$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 176
$ws.Cells.Item(176, 1).Value = 175
$ws.Cells.Item(176, 2).Value = 1
$ws.Cells.Item(176, 3).Value = "2024-06-18 16:16:16"
$ws.Cells.Item(176, 4).Value = 200
$ws.Cells.Item(176, 5).Value = 22

# Row 177
$ws.Cells.Item(177, 1).Value = 176
$ws.Cells.Item(177, 2).Value = 2
$ws.Cells.Item(177, 3).Value = "2024-06-18 16:16:16"
$ws.Cells.Item(177, 4).Value = 200
$ws.Cells.Item(177, 5).Value = 0
